$d = $word.ActiveDocument

# 1) Update the date in the header line
$d.Content.Find.Execute(
    "17.08.24",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "16.08.24", 2) | Out-Null

# 2) Update the paper title
$d.Content.Find.Execute(
    "Faster Machine Unlearning via Natural Gradient Descent",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "On the Geometry of Deep Learning", 2) | Out-Null

# 3) Remove the old body paragraphs 3 through 12 (the whole unlearning discussion
#    after the title, up to -- but not including -- the arxiv link paragraph).
$startPara = $d.Paragraphs.Item(3)
$endPara = $d.Paragraphs.Item(12)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# 4) Insert 5 brand-new empty paragraphs right before the (now-third) paragraph,
#    which is the arxiv link paragraph. Inserting fresh paragraphs here (rather
#    than overwriting the text of the old paragraph 3, whose run happened to
#    carry an xml:space="preserve" run attribute) keeps the new runs' whitespace
#    handling correct/minimal.
$newTexts = @(
    "אני ממש אוהב מאמרים שחוקרים מה שקורה בתוך המודלים העמוקים שלנו - הרי לדעתי זה התנאי הכרחי לכך שנוכל להתחיל באמת לסמוך על- AI (לפחות חלקית). ואכן הכותבים מדגישים כי למידה עמוקה, על אף הישגיה המרשימים במגוון תחומים, נשארת עדיין בגדר ""קופסה שחורה"" עם הבנה חלקית בלבד של אופן פעולתה.",
    "המחברים מנסים להסביר מודלים עמוקים באמצעות ספליינים אפיניים (Affine Splines) שהן למעשה פונקציות רציפות ולינאריות למקוטעין במרחב רב מימד. המחקר מתבונן ברשתות נוירונים מזווית גיאומטרית באמצעות ניתוח של חלוקות הנוצרות על ידי ספליינים אפיניים, המקרבות אותן (הרשתות). ",
    "בפרט המחברים דנים בחלוקות של מרחב הקלט לפי הקטגוריות שלו הנוצרות על ידי ייצוג לטנטי (השכבה האחרונה לפני שכבת הסיווג) של הרשת. הבנת החלוקה הזו מסייעת להסביר כיצד רשתות עמוקות לומדות ומייצרות חיזוים עבור קלטים שונים. ",
    "המחברים גם דנים במבנים גיאומטריים הנוצרים על ידי משקלי המודל במרחב הלוס (כלומר מנתחים את פונקציית הלוס למשקלי הרשת השונים). בנוסף המאמר גם מדבר על החלוקות הנוצרות במרחב משקולות המודל בשכבות שונות לאתחולי רשת שונים וגם לאימון עם ובלי BatchNorm. כמובן שזה נעשה על דוגמאות מלאכותיות(toy examples) בעלי מימד נמוך. ויש עוד מספר ניתוחים גיאומטרים די מעניינים במאמר.",
    "מעניין כי המחברים כותבים כי אחת המטרות המרכזיות של המחקר היא לדרבן מתמטיקאים לעסוק בניתוח גיאומטרי של רשתות עמוקות."
)

foreach ($t in $newTexts) {
    $urlPara = $d.Paragraphs.Item(3)
    $insertPoint = $d.Range($urlPara.Range.Start, $urlPara.Range.Start)
    $insertPoint.InsertParagraphBefore()
}

for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $d.Paragraphs.Item(3 + $i).Range.Text = $newTexts[$i]
}

# 5) Update the arxiv link (now the last paragraph)
$d.Content.Find.Execute(
    "https://arxiv.org/abs/2407.08169",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://arxiv.org/abs/2408.04809", 2) | Out-Null

Write-Output ("Final paragraph count=" + $d.Paragraphs.Count)
